$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start with a clean sheet so the shared-string table is rebuilt in the
# natural row-by-row order that matches the target layout.
$ws.Cells.Clear()

# Full data set (A = original_value, B = real_value_english) for rows 1..25
$data = @(
    @("original_value", "real_value_english"),
    @("utilitiesCommuncation", "utilitiesCommunication"),
    @("Environment", "environment"),
    @("biota", "biota"),
    @("boundaries", "boundaries"),
    @("climatologyMeteorologyAtmosphere", "climatologyMeteorologyAtmosphere"),
    @("economy", "economy"),
    @("elevation", "elevation"),
    @("environment", "environment"),
    @("farming", "farming"),
    @("geoscientificInformation", "geoscientificInformation"),
    @("health", "health"),
    @("imageryBaseMapsEarthCover", "imageryBaseMapsEarthCover"),
    @("intelligenceMilitary", "intelligenceMilitary"),
    @("inlandWaters", "inlandWaters"),
    @("location", "location"),
    @("oceans", "oceans"),
    @("planningCadastre", "planningCadastre"),
    @("society", "society"),
    @("structure", "structure"),
    @("transportation", "transportation"),
    @("utilitiesCommunication", "utilitiesCommunication"),
    @("inlandwaters", "inlandWaters"),
    @(" boundaries", "boundaries"),
    @(" inlandwaters", "inlandWaters")
)

$row = 1
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Column B width change (stored OOXML width of 26 characters).
# Excel's ColumnWidth property value and the stored "width" attribute in the
# XML differ by the standard column-width padding correction, so the input
# value is adjusted so the saved file ends up with width="26".
$ws.Columns.Item(2).ColumnWidth = 25.14

# Sheet view changes: remove topLeftCell freeze-scroll position and change selection
$ws.Range("A4:B25").Select()
